# gas biasing circuit - add total/1차/2차/3차/현재 tracking columns (I:M)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Header row (row 1), columns I:M ----
$ws.Range("I1").Value = "total"
$ws.Range("J1").Value = "1차"
$ws.Range("K1").Value = "2차"
$ws.Range("L1").Value = "3차"
$ws.Range("M1").Value = "현재"

# ---- "total" (I) column starting values, per row ----
$iValues = @{
    2 = 10000;  3 = 1000;   4 = 10000;  5 = 10000;  6 = 5000;
    7 = 500;    8 = 500;    9 = 60;     10 = 500;   11 = 200;
    12 = 100;   13 = 200;   14 = 500;   15 = 500;   16 = 500;
    17 = 500;   18 = 500;   19 = 10000; 20 = 10000; 21 = 10000;
    22 = 10000; 23 = 10000; 24 = 10000; 25 = 10000; 26 = 10000;
    27 = 10000; 28 = 10000; 29 = 10000; 30 = 100;   31 = 200
}

for ($r = 2; $r -le 31; $r++) {
    $ws.Cells.Item($r, 9).Value = $iValues[$r]

    # J (1차), K (2차), L (3차) = 20 * Qty, except row 13 (no 1차 usage) and
    # row 11 whose 2차 usage (K) was hand-overridden to a literal count.
    if ($r -ne 13) {
        $ws.Cells.Item($r, 10).Formula = "=20*C$r"
    }
    if ($r -eq 11) {
        $ws.Cells.Item($r, 11).Value = 84
    } else {
        $ws.Cells.Item($r, 11).Formula = "=20*C$r"
    }
    $ws.Cells.Item($r, 12).Formula = "=20*C$r"

    # M (현재) = total - 1차 - 2차 - 3차
    $ws.Cells.Item($r, 13).Formula = "=I$r-J$r-K$r-L$r"
}

# ---- Apply the built-in "good" / "bad" cell styles ----
# Row 11 stands out as a shortage (usage exceeded the tracked total) -> Bad
$ws.Range("A11:M11").Style = "나쁨"

# All other populated rows -> Good, except rows 4, 5 and 20 which were left
# with the default style (row 20 only picked up the style on its I cell).
$goodRows = @(2,3,6,7,8,9,10,12,13,14,15,16,17,18,19,21,22,23,24,25,26,27,28,29,30,31)
foreach ($r in $goodRows) {
    $ws.Range("A" + $r + ":M" + $r).Style = "좋음"
}
$ws.Cells.Item(20, 9).Style = "좋음"

# ---- Selection, matching the saved session ----
$ws.Range("G40").Select()
